# "starting reformulation with planninghorizon set"
# Adds a new "parameters" worksheet (planning_horizon / time_delta_minutes)
# after the existing node/edge/device sheets, and makes it the active sheet.

$wb = $excel.ActiveWorkbook

# --- add the new sheet after the last existing sheet ------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "parameters"

# --- header row ---------------------------------------------------------
$ws.Range("A1").Value = "parameter"
$ws.Range("B1").Value = "value"
$ws.Range("C1").Value = "comment"

# --- data rows (column order chosen to match authoring / shared-string
#     allocation order: A2, C3, C2, A3, B2, B3) -----------------------
$ws.Range("A2").Value = "planning_horizon"
$ws.Range("C3").Value = "minutes"
$ws.Range("C2").Value = "number of time steps in planning horizon"
$ws.Range("A3").Value = "time_delta_minutes"
$ws.Range("B2").Value = 48
$ws.Range("B3").Value = 15

# --- column widths (closest values reachable through the ColumnWidth
#     property, which is quantized to 1/6 character increments) --------
$ws.Columns.Item(1).ColumnWidth = 14.1667
$ws.Columns.Item(3).ColumnWidth = 33.6667

# --- selection / activation so tabSelected + activeTab move onto the
#     new sheet (and off of the previously-selected "node" sheet) ------
$ws.Range("A6").Select() | Out-Null
$ws.Activate() | Out-Null
